# Updates the cryptocurrency price/volume table (GitHub Actions style refresh).
# Columns D (Price) and E (Volume 1h) hold values that look numeric but must
# stay plain text (as in the original inline strings). For values that are
# valid numeric literals, a leading apostrophe forces Excel to keep them as
# text; the style is then reset to "Normal" so no extra number-format / quote
# -prefix style gets attached to the cell (matching the original, unstyled
# cells exactly). Values that already aren't valid numbers (e.g. "30.583.77"
# with two dots, or the padded "  +1.30%  " percentages) are simply assigned.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "30.583.77"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "1.924.55"
$ws.Range("E3").Value = "  +3.73%  "
$ws.Range("E4").Value = "  -0.02%  "
$ws.Range("D5").Value = "'247.44"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +5.00%  "
$ws.Range("E6").Value = "  -0.11%  "
$ws.Range("D7").Value = "'0.4736"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.56%  "
$ws.Range("D8").Value = "'0.2908"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +3.21%  "
$ws.Range("D9").Value = "'0.06779"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +5.87%  "
$ws.Range("D10").Value = "'105.00"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +8.44%  "
$ws.Range("D11").Value = "'18.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +1.05%  "
$ws.Range("D12").Value = "1.910.35"
$ws.Range("E12").Value = "  +3.03%  "
$ws.Range("D13").Value = "'0.07726"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +2.51%  "
$ws.Range("E14").Value = "  +6.69%  "
$ws.Range("D15").Value = "'0.6724"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +5.18%  "
$ws.Range("D16").Value = "'287.15"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.39%  "
$ws.Range("D17").Value = "30.622.18"
$ws.Range("E17").Value = "  +1.43%  "
$ws.Range("D18").Value = "'0.000007630"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +3.31%  "
$ws.Range("D19").Value = "'12.98"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +2.00%  "
$ws.Range("D20").Value = "'1.000"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -0.06%  "
$ws.Range("D21").Value = "2.160.19"
$ws.Range("E21").Value = "  +2.69%  "
$ws.Range("E22").Value = "  +8.93%  "
$ws.Range("E23").Value = "  +0.00%  "
$ws.Range("D24").Value = "'6.313"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +4.41%  "
$ws.Range("D25").Value = "'9.396"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  +3.39%  "
$ws.Range("D26").Value = "'168.14"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.89%  "
$ws.Range("D27").Value = "'20.76"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +7.08%  "
$ws.Range("D28").Value = "'2.153"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +11.52%  "
$ws.Range("D29").Value = "'0.1082"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.04%  "
$ws.Range("D30").Value = "'1.367"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +3.10%  "
$ws.Range("D31").Value = "'4.211"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.72%  "
$ws.Range("D32").Value = "'4.140"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +8.84%  "
$ws.Range("D33").Value = "'0.05057"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  +2.89%  "
$ws.Range("B34").Value = "ImmutableX"
$ws.Range("C34").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D34").Value = "'0.7433"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  +2.39%  "
$ws.Range("B35").Value = "ARBITRUM"
$ws.Range("C35").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D35").Value = "'1.163"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +4.26%  "
$ws.Range("D36").Value = "'0.02080"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +8.25%  "
$ws.Range("D37").Value = "'2.746"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("D38").Value = "'2.692"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +0.90%  "
$ws.Range("D39").Value = "'2.070"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  +5.37%  "
$ws.Range("D40").Value = "'111.32"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +5.50%  "
$ws.Range("D41").Value = "'0.8821"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  +1.96%  "
$ws.Range("D42").Value = "'5.966"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +6.08%  "
$ws.Range("D43").Value = "'0.4369"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +7.65%  "
$ws.Range("E44").Value = "  -0.22%  "
$ws.Range("D45").Value = "'67.26"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +2.77%  "
$ws.Range("D46").Value = "'7.275"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  +2.96%  "
$ws.Range("D47").Value = "'9.364"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +4.27%  "
$ws.Range("D48").Value = "'48.41"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +18.27%  "
$ws.Range("E49").Value = "  +3.97%  "
$ws.Range("D50").Value = "'35.22"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.48%  "
$ws.Range("D51").Value = "'0.4044"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +8.29%  "
